$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1764705882352941
$ws.Range("C2").Value = 0.6470588235294118
$ws.Range("P2").Value = 0.1764705882352941
$ws.Range("P3").Value = 0.5454545454545454
$ws.Range("S3").Value = 0.4545454545454545
$ws.Range("B6").Value = 0.1333333333333333
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.06666666666666667
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.4666666666666667
$ws.Range("B7").Value = 0.03846153846153846
$ws.Range("F7").Value = 0.03846153846153846
$ws.Range("J7").Value = 0.1153846153846154
$ws.Range("O7").Value = 0.03846153846153846
$ws.Range("Q7").Value = 0.1538461538461539
$ws.Range("R7").Value = 0.07692307692307693
$ws.Range("S7").Value = 0.5384615384615384
$ws.Range("B8").Value = 0.1052631578947368
$ws.Range("D8").Value = 0.02631578947368421
$ws.Range("F8").Value = 0.1052631578947368
$ws.Range("J8").Value = 0.1052631578947368
$ws.Range("O8").Value = 0.02631578947368421
$ws.Range("Q8").Value = 0.1052631578947368
$ws.Range("R8").Value = 0.05263157894736842
$ws.Range("S8").Value = 0.4736842105263158
$ws.Range("D9").Value = 0.08333333333333333
$ws.Range("F9").Value = 0.08333333333333333
$ws.Range("J9").Value = 0.08333333333333333
$ws.Range("O9").Value = 0.08333333333333333
$ws.Range("Q9").Value = 0.25
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.3333333333333333
$ws.Range("B10").Value = 0.08974358974358974
$ws.Range("F10").Value = 0.07692307692307693
$ws.Range("J10").Value = 0.1538461538461539
$ws.Range("R10").Value = 0.05128205128205128
$ws.Range("S10").Value = 0.4615384615384616
$ws.Range("G11").Value = 0.3095238095238095
$ws.Range("J11").Value = 0.07142857142857142
$ws.Range("K11").Value = 0.3095238095238095
$ws.Range("L11").Value = 0.3095238095238095
$ws.Range("G12").Value = 0.9230769230769231
$ws.Range("J12").Value = 0.07692307692307693
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("S13").Value = 0.3333333333333333
$ws.Range("H15").Value = 0.3076923076923077
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.1538461538461539
$ws.Range("O15").Value = 0.1538461538461539
$ws.Range("S15").Value = 0.3076923076923077
$ws.Range("H16").Value = 0.2727272727272727
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.1818181818181818
$ws.Range("K16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.3636363636363636
$ws.Range("H17").Value = 0.3703703703703703
$ws.Range("I17").Value = 0.03703703703703703
$ws.Range("J17").Value = 0.2592592592592592
$ws.Range("K17").Value = 0.2222222222222222
$ws.Range("S17").Value = 0.1111111111111111
$ws.Range("H18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.5555555555555556
$ws.Range("O18").Value = 0.1111111111111111
$ws.Range("S18").Value = 0.2222222222222222
$ws.Range("H19").Value = 0.1666666666666667
$ws.Range("I19").Value = 0.075
$ws.Range("J19").Value = 0.325
$ws.Range("K19").Value = 0.1833333333333333
$ws.Range("M19").Value = 0.025
$ws.Range("O19").Value = 0.04166666666666666
$ws.Range("S19").Value = 0.1833333333333333
